$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40
$prevRow = $row - 1

# Copy the formatting from the previous row so the new row matches
# the existing style (date number format in column A, etc.)
$ws.Range("A$prevRow`:E$prevRow").Copy()
$ws.Range("A$row`:E$row").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 45944
$ws.Cells.Item($row, 2).Value = "21,5965"
$ws.Cells.Item($row, 3).Value = "15,5224"
$ws.Cells.Item($row, 4).Value = "15,3715"
$ws.Cells.Item($row, 5).Value = "15,3715"
